$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first 10 data rows (years 2000-2009), which are in rows 2 to 11.
# This shifts the remaining rows (2010-2019, originally rows 12-21) up to rows 2-11.
$ws.Range("A2:B11").EntireRow.Delete()

$wb.Save()
